$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 767.8823
$ws.Range("I28").Value = 758.0769
$ws.Range("J28").Value = 799.75
$ws.Range("K28").Value = 758.0769
$ws.Range("L28").Value = 799.75
$ws.Range("M28").Value = -273.0769
$ws.Range("N28").Value = -1769.75

$ws.Range("H87").Value = 16984.38
$ws.Range("J87").Value = 16984.38
$ws.Range("L87").Value = 16984.38
$ws.Range("N87").Value = -19480.38

$ws.Range("H90").Value = 16984.38
$ws.Range("J90").Value = 16984.38
$ws.Range("L90").Value = 50953.14
$ws.Range("N90").Value = -63433.14

$ws.Range("H98").Value = 1412.5625
$ws.Range("I98").Value = 1420.5
$ws.Range("J98").Value = 1399.3334
$ws.Range("K98").Value = 1420.5
$ws.Range("L98").Value = 1399.3334
$ws.Range("M98").Value = 77.5
$ws.Range("N98").Value = -4395.3334

$ws.Range("H106").Value = 837.1429000000001
$ws.Range("I106").Value = 768.8889
$ws.Range("J106").Value = 960
$ws.Range("K106").Value = 768.8889
$ws.Range("L106").Value = 960
$ws.Range("M106").Value = -137.8889
$ws.Range("N106").Value = -2222

$ws.Range("H112").Value = 1059.5306
$ws.Range("J112").Value = 1087.711
$ws.Range("L112").Value = 3263.133
$ws.Range("N112").Value = -5479.133

$ws.Range("H122").Value = 1412.5625
$ws.Range("I122").Value = 1420.5
$ws.Range("J122").Value = 1399.3334
$ws.Range("K122").Value = 4261.5
$ws.Range("L122").Value = 4198.0002
$ws.Range("M122").Value = -1811.5
$ws.Range("N122").Value = -9098.0002

$ws.Range("H125").Value = 1682997.9
$ws.Range("I125").Value = 3260
$ws.Range("J125").Value = 2942801.2
$ws.Range("K125").Value = 29340
$ws.Range("L125").Value = 26485210.8
$ws.Range("M125").Value = -26880
$ws.Range("N125").Value = -26490130.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1246.16
$ws.Range("I45").Value = 1206.4286
$ws.Range("K45").Value = 1206.4286
$ws.Range("M45").Value = -829.4286

$ws.Range("H74").Value = 7481.05
$ws.Range("I74").Value = 10065.923
$ws.Range("J74").Value = 2680.5715
$ws.Range("K74").Value = 10065.923
$ws.Range("L74").Value = 2680.5715
$ws.Range("M74").Value = -9191.923000000001
$ws.Range("N74").Value = -4428.5715

$ws.Range("H77").Value = 7481.05
$ws.Range("I77").Value = 10065.923
$ws.Range("J77").Value = 2680.5715
$ws.Range("K77").Value = 50329.61500000001
$ws.Range("L77").Value = 13402.8575
$ws.Range("M77").Value = -45961.61500000001
$ws.Range("N77").Value = -22138.8575

$ws.Range("H110").Value = 102282.2
$ws.Range("I110").Value = 167870.33
$ws.Range("J110").Value = 3900
$ws.Range("K110").Value = 167870.33
$ws.Range("L110").Value = 3900
$ws.Range("M110").Value = -165825.33
$ws.Range("N110").Value = -7990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2793.1667
$ws.Range("I105").Value = 3072.25
$ws.Range("J105").Value = 2235
$ws.Range("K105").Value = 3072.25
$ws.Range("L105").Value = 2235
$ws.Range("M105").Value = -1325.25
$ws.Range("N105").Value = -5729

$ws.Range("H134").Value = 1195.9814
$ws.Range("I134").Value = 874.5349
$ws.Range("J134").Value = 2452.5454
$ws.Range("K134").Value = 2623.6047
$ws.Range("L134").Value = 7357.6362
$ws.Range("M134").Value = -88.60469999999987
$ws.Range("N134").Value = -12427.6362

$ws.Range("H140").Value = 41883.5
$ws.Range("J140").Value = 41883.5
$ws.Range("L140").Value = 41883.5
$ws.Range("N140").Value = -52243.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2116.394
$ws.Range("I58").Value = 1492.6316
$ws.Range("J58").Value = 2962.9285
$ws.Range("K58").Value = 1492.6316
$ws.Range("L58").Value = 2962.9285
$ws.Range("M58").Value = -1289.6316
$ws.Range("N58").Value = -3368.9285

$ws.Range("H99").Value = 1784465
$ws.Range("I99").Value = 2001210.5
$ws.Range("J99").Value = 50500
$ws.Range("K99").Value = 2001210.5
$ws.Range("L99").Value = 50500
$ws.Range("M99").Value = -1999712.5
$ws.Range("N99").Value = -53496

$ws.Range("H126").Value = 1784465
$ws.Range("I126").Value = 2001210.5
$ws.Range("J126").Value = 50500
$ws.Range("K126").Value = 6003631.5
$ws.Range("L126").Value = 151500
$ws.Range("M126").Value = -6001161.5
$ws.Range("N126").Value = -156440

$ws.Range("H132").Value = 2982.6667
$ws.Range("I132").Value = 2316.0625
$ws.Range("J132").Value = 3952.2727
$ws.Range("K132").Value = 6948.1875
$ws.Range("L132").Value = 11856.8181
$ws.Range("M132").Value = -4418.1875
$ws.Range("N132").Value = -16916.8181

$ws.Range("H136").Value = 2116.394
$ws.Range("I136").Value = 1492.6316
$ws.Range("J136").Value = 2962.9285
$ws.Range("K136").Value = 4477.8948
$ws.Range("L136").Value = 8888.7855
$ws.Range("M136").Value = -1927.8948
$ws.Range("N136").Value = -13988.7855

$ws.Range("H138").Value = 35787.5
$ws.Range("J138").Value = 35787.5
$ws.Range("L138").Value = 35787.5
$ws.Range("N138").Value = -46067.5

$ws.Range("H140").Value = 56173.168
$ws.Range("J140").Value = 56173.168
$ws.Range("L140").Value = 56173.168
$ws.Range("N140").Value = -66533.16800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 975707.75
$ws.Range("I5").Value = 511.93332
$ws.Range("J5").Value = 5851686.5
$ws.Range("K5").Value = 1535.79996
$ws.Range("L5").Value = 17555059.5
$ws.Range("M5").Value = -1423.79996
$ws.Range("N5").Value = -17555283.5

$ws.Range("H135").Value = 975707.75
$ws.Range("I135").Value = 511.93332
$ws.Range("J135").Value = 5851686.5
$ws.Range("K135").Value = 4607.39988
$ws.Range("L135").Value = 52665178.5
$ws.Range("M135").Value = -2072.39988
$ws.Range("N135").Value = -52670248.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 35001.5
$ws.Range("I4").Value = 20003
$ws.Range("K4").Value = 20003
$ws.Range("M4").Value = -19891

$ws.Range("H132").Value = 1846.7368
$ws.Range("I132").Value = 1405.9333
$ws.Range("J132").Value = 3499.75
$ws.Range("K132").Value = 4217.7999
$ws.Range("L132").Value = 10499.25
$ws.Range("M132").Value = -1687.7999
$ws.Range("N132").Value = -15559.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 102015.1
$ws.Range("I46").Value = 169333.5
$ws.Range("J46").Value = 1037.5
$ws.Range("K46").Value = 169333.5
$ws.Range("L46").Value = 1037.5
$ws.Range("M46").Value = -169145.5
$ws.Range("N46").Value = -1413.5

$ws.Range("H132").Value = 6066.0923
$ws.Range("I132").Value = 7890.946
$ws.Range("K132").Value = 23672.838
$ws.Range("M132").Value = -21142.838

$ws.Range("H136").Value = 7777590
$ws.Range("I136").Value = 45824.39
$ws.Range("J136").Value = 16669120
$ws.Range("K136").Value = 137473.17
$ws.Range("L136").Value = 50007360
$ws.Range("M136").Value = -134923.17
$ws.Range("N136").Value = -50012460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 41425.8
$ws.Range("I122").Value = 51172.35
$ws.Range("J122").Value = 2439.6
$ws.Range("K122").Value = 153517.05
$ws.Range("L122").Value = 7318.799999999999
$ws.Range("M122").Value = -151067.05
$ws.Range("N122").Value = -12218.8
